$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values must be forced to text storage (to match the original
# inline-string cells) since many look like plain numbers (e.g. "246.05",
# "0.660", "12.40") and would otherwise be auto-coerced by Excel into
# numeric cells, losing formatting like trailing/leading zeros. We force
# text via NumberFormat "@" then restore the default (unstyled) look with
# ClearFormats() so no stray style index is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.199.68"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.06"
$ws.Range("D3").ClearFormats()

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.660"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.73%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.25"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.346"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +13.12%  "

$ws.Range("E11").Value = "  +2.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0992"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.167.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.56%  "

$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.892.62"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.82"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.183.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0818"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "240.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.58%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +26.44%  "

$ws.Range("E26").Value = "  +1.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.43"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.180.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +22.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.947"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +14.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0560"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.08"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("E37").Value = "  -2.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.02"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("E39").Value = "  +1.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.10"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0208"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.13"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0629"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "89.48"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.333.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.61"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +38.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("E49").Value = "  +1.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.47"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.078.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.39%  "
